$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7180846140156819
$ws.Range("C2").Value = 0.6318612555126804
$ws.Range("D2").Value = 0.5993536909734288

$ws.Range("B3").Value = 0.7570384412011173
$ws.Range("C3").Value = -0.7189261932389713
$ws.Range("D3").Value = -0.6674135141025708

$ws.Range("B4").Value = 0.8865791316082355
$ws.Range("C4").Value = 0.5217357206468962
$ws.Range("D4").Value = 0.5849859299172766

$ws.Range("B5").Value = -0.8002462272016042
$ws.Range("C5").Value = -0.5471126479043555
$ws.Range("D5").Value = -0.5887247801971668

$ws.Range("B6").Value = 0.8208186911390987
$ws.Range("C6").Value = -0.5761031271683708
$ws.Range("D6").Value = -0.5741786760964228

$ws.Range("B7").Value = 0.7876702562287906
$ws.Range("C7").Value = -0.7594794030903435
$ws.Range("D7").Value = -0.600815490967664

$ws.Range("B8").Value = 0.6777699892607838
$ws.Range("C8").Value = -0.6067459581690833

$ws.Range("B9").Value = -0.6676936083905627
$ws.Range("C9").Value = -0.6138190161125965
$ws.Range("D9").Value = 0.5423844290542225
